$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point rounding on the existing row 80 timestamp
$ws.Range("A80").Value = 44393.77293260532

# Append the new row 81 with the latest retrieved data
$ws.Range("A81").Value = 44394.77242421631
$ws.Range("B81").Value = 80252
$ws.Range("C81").Value = 67689
$ws.Range("D81").Value = 3672
$ws.Range("E81").Value = 2205
$ws.Range("F81").Value = 1588
$ws.Range("G81").Value = 21101
$ws.Range("H81").Value = 1614
$ws.Range("I81").Value = 902
$ws.Range("J81").Value = 200

# Match the date-formatted style used by the rest of column A
$ws.Range("A81").NumberFormat = $ws.Range("A80").NumberFormat
